# ozon fixes 30.10.2025 part 1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update order / sum values for rows 2-4
$ws.Range("A2").Value = 2056204
$ws.Range("B2").Value = 18457

$ws.Range("A3").Value = 2059046
$ws.Range("B3").Value = 19292

$ws.Range("A4").Value = 2083871
$ws.Range("B4").Value = 25646

# Rows 5 and 6 no longer have data - clear them out entirely
$ws.Range("A5:E5").ClearContents()
$ws.Range("A6:E6").ClearContents()

# Move the active selection to B9
$ws.Range("B9").Select()
